$wb = $excel.ActiveWorkbook

# Insert the new "Asset" sheet right before "LoginPage" so the final tab
# order is Contract, Asset, LoginPage - regardless of which sheet happens
# to be active when this script runs.
$loginPage = $wb.Worksheets.Item("LoginPage")
$assetSheet = $wb.Worksheets.Add($loginPage)
$assetSheet.Name = "Asset"

# Header row
$assetSheet.Range("A1").Value = "AssetName"
$assetSheet.Range("B1").Value = "AccountName"
$assetSheet.Range("C1").Value = "ContactName"

# Data row
$assetSheet.Range("A2").Value = "test asset name"
$assetSheet.Range("B2").Value = "test"
$assetSheet.Range("C2").Value = "test"

# Column B reuses the existing accounting-style number format (style index 1
# in the workbook), matching the "AccountName"/amount column on the
# Contract sheet.
$assetSheet.Range("B1:B2").NumberFormat = "#,##0_);(#,##0)"

# Update the selection on the Contract sheet (A3:XFD3 -> A1:B2).
$contract = $wb.Worksheets.Item("Contract")
[void]$contract.Range("A1:B2").Select()

# Leave the new Asset sheet active, with B1 selected.
[void]$assetSheet.Select()
[void]$assetSheet.Range("B1").Select()
